$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Journey time" row (currently row 61) to hold the
# second ("double") meal entry for the day -> new row 61 = MEAL/Lunch row.
$ws.Rows.Item(61).Insert()

# Values for the new MEAL row (B=label, C..G = L/blank/L/L/L)
$ws.Range("B61").Value = "MEAL"
$ws.Range("C61").Value = "L"
$ws.Range("E61").Value = "L"
$ws.Range("F61").Value = "L"
$ws.Range("G61").Value = "L"

Write-Output "done"
